$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value, taken from the commit diff.
# Values are plain-text quote (price / percentage) columns (D = Price, E = Volume(1h)),
# stored in the sheet as inline/shared text, not numbers -- so each value is written with a
# leading apostrophe to force text entry (avoids Excel auto-converting "285.42" / "3.10%" into
# a float), and then the cell style is reset to "Normal" so the quote-prefix flag added by the
# apostrophe does not leave a stray style on the cell.
$updates = [ordered]@{
    "D2" = "285.42"
    "E2" = "3.10%"
    "E3" = "4.53%"
    "D4" = "5.007"
    "E4" = "2.76%"
    "D5" = "0.06488"
    "E5" = "1.16%"
    "D6" = "7.234"
    "E6" = "4.21%"
    "D7" = "1.343"
    "E7" = "13.61%"
    "D8" = "0.9133"
    "E8" = "4.25%"
    "D9" = "0.1542"
    "E9" = "0.54%"
    "D10" = "0.06467"
    "E10" = "25.34%"
    "D11" = "0.07617"
    "E11" = "1.90%"
    "D12" = "0.02983"
    "E12" = "1.06%"
    "D13" = "0.08971"
    "E13" = "-0.10%"
    "D14" = "0.001596"
    "E14" = "1.60%"
    "D15" = "0.0006560"
    "E15" = "3.15%"
    "D16" = "0.006036"
    "E16" = "-1.44%"
    "D17" = "3.461"
    "E17" = "-0.57%"
    "D18" = "3.366"
    "E18" = "1.87%"
    "E19" = "-1.44%"
    "D21" = "0.1341"
    "E21" = "1.57%"
    "D22" = "3.973"
    "E22" = "1.68%"
    "E23" = "3.63%"
    "D24" = "0.04461"
    "E24" = "1.35%"
    "D25" = "0.001182"
    "E25" = "0.42%"
    "D26" = "0.004328"
    "E26" = "11.95%"
    "E28" = "-9.23%"
    "D29" = "0.0001636"
    "E29" = "-15.71%"
    "D40" = "0.04149"
    "E40" = "-0.43%"
    "D41" = "0.006732"
    "E41" = "-1.03%"
    "D42" = "0.1234"
    "E42" = "5.21%"
    "D43" = "0.002109"
    "E43" = "3.39%"
    "D44" = "0.01181"
    "E44" = "3.07%"
    "D45" = "0.00005377"
    "E45" = "1.42%"
    "E46" = "-0.03%"
    "E47" = "21.16%"
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $newValue
    $range.Style = "Normal"
}

